# Apply updated cryptocurrency price/volume data (auto-generated)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.238.38"
$ws.Range("E2").Value = "  +5.40%  "

$ws.Range("D3").Value = "2.792.46"
$ws.Range("E3").Value = "  +5.82%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "117.15"
$ws.Range("E5").Value = "  +4.69%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "341.29"
$ws.Range("E6").Value = "  +4.72%  "

$ws.Range("E7").Value = "  +5.45%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("E9").Value = "  +5.96%  "

$ws.Range("E10").Value = "  +6.71%  "

$ws.Range("E11").Value = "  +7.50%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.12"
$ws.Range("E12").Value = "  +0.11%  "

$ws.Range("E13").Value = "  +2.45%  "

$ws.Range("E14").Value = "  +1.69%  "

$ws.Range("D15").Value = "3.235.24"
$ws.Range("E15").Value = "  +6.03%  "

$ws.Range("D16").Value = "2.778.32"
$ws.Range("E16").Value = "  +5.25%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.887"
$ws.Range("E17").Value = "  +3.94%  "

$ws.Range("D18").Value = "52.077.88"
$ws.Range("E18").Value = "  +5.18%  "

$ws.Range("E19").Value = "  +11.57%  "

$ws.Range("E20").Value = "  +1.98%  "

$ws.Range("E21").Value = "  +4.23%  "

$ws.Range("D22").Value = "0.0₃0988"
$ws.Range("E22").Value = "  +4.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "278.45"
$ws.Range("E23").Value = "  +3.64%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.36"
$ws.Range("E24").Value = "  +1.80%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.82"
$ws.Range("E25").Value = "  +10.32%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.93"
$ws.Range("E26").Value = "  +3.16%  "

$ws.Range("E27").Value = "  -0.08%  "

$ws.Range("E28").Value = "  +0.10%  "

$ws.Range("E29").Value = "  +1.17%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.142"
$ws.Range("E30").Value = "  +2.55%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.00"
$ws.Range("E31").Value = "  +0.90%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.41"
$ws.Range("E32").Value = "  +1.72%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.71"
$ws.Range("E33").Value = "  +4.02%  "

$ws.Range("E34").Value = "  +2.22%  "

$ws.Range("E35").Value = "  +4.90%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  -0.01%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.03"
$ws.Range("E37").Value = "  -0.06%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.98"
$ws.Range("E38").Value = "  +0.72%  "

$ws.Range("E39").Value = "  +5.67%  "

$ws.Range("E40").Value = "  +29.09%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0372"
$ws.Range("E41").Value = "  +12.96%  "

$ws.Range("E42").Value = "  +5.25%  "

$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "23.44"
$ws.Range("E43").Value = "  +2.18%  "

$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.116"
$ws.Range("E44").Value = "  +4.33%  "

$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "127.13"
$ws.Range("E45").Value = "  -1.01%  "

$ws.Range("D46").Value = "2.107.10"
$ws.Range("E46").Value = "  +2.42%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.35"
$ws.Range("E47").Value = "  +2.63%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.56"
$ws.Range("E49").Value = "  +6.72%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.924"
$ws.Range("E50").Value = "  +22.94%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.96"
$ws.Range("E51").Value = "  +1.19%  "
